$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = " CKO oQf"
$ws.Range("A3").Value = "Qt"
$ws.Range("A4").Value = "9 fO2v 2x D"
$ws.Range("A5").Value = "C P2IWmhi"
$ws.Range("A6").Value = "9l"
